{"js": "// Bold + color (\"#2C3E50\") the quantitative/impact metrics inside specific\n// bullet paragraphs of the resume, matching the \"hybrid bold + color\n// highlighting\" treatment described in the commit message (percentages,\n// dollar amounts, etc.). Only the exact paragraphs touched by the diff are\n// modified; other occurrences of the same numbers elsewhere in the document\n// (e.g. the Professional Summary) are left untouched because the search is\n// scoped to each individual paragraph.\n\nconst HIGHLIGHT_COLOR = \"#2C3E50\";\n\n// Each entry: a snippet of text that uniquely identifies the target\n// paragraph (its leading text), plus the ordered list of metric substrings\n// inside that paragraph that must become bold + colored.\nconst targets = [\n  {\n    lead: \"\u2022 Discovered systematic race coding errors\",\n    metrics: [\"23%\", \"64%\"],\n  },\n  {\n    lead: \"\u2022 Utilized advanced sampling methods\",\n    metrics: [\"\\u00B14.2%\", \"\\u00B12.1%\", \"71%\", \"87%\"],\n  },\n  {\n    lead: \"\u2022 Trigonometric algorithm for boundary estimation\",\n    metrics: [\"73.5%\", \"$4.7M\"],\n  },\n  {\n    lead: \"\u2022 Built real-time FEC analysis systems\",\n    metrics: [\"$2\"],\n  },\n  {\n    lead: \"\u2022 Algorithmic innovation: Pioneered trigonometric\",\n    metrics: [\"73.5%\"],\n  },\n  {\n    lead: \"\u2022 $4.7M savings enabled\",\n    metrics: [\"$4.7M\"],\n  },\n  {\n    lead: \"\u2022 178% accuracy improvement\",\n    metrics: [\"178%\"],\n  },\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Collect the search-result ranges for every metric in every paragraph\n// first, then format them all in a second pass (fewer sync round-trips).\nconst foundRanges = [];\n\nfor (const target of targets) {\n  const paragraph = paragraphs.items.find((p) => p.text.indexOf(target.lead) === 0);\n  if (!paragraph) {\n    continue;\n  }\n  for (const metric of target.metrics) {\n    const results = paragraph.search(metric, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    foundRanges.push(results);\n  }\n}\n\nawait context.sync();\n\nfor (const results of foundRanges) {\n  if (results.items.length > 0) {\n    results.items[0].font.set({ bold: true, color: HIGHLIGHT_COLOR });\n  }\n}\n\nawait context.sync();\n", "ps1": "# Bold + color (\"2C3E50\") the quantitative/impact metrics inside specific\n# bullet paragraphs of the resume, matching the \"hybrid bold + color\n# highlighting\" treatment described in the commit message.\n#\n# Strategy: for each target paragraph (identified by its distinctive\n# leading text so the script is robust to paragraph renumbering), run a\n# Find scoped to that paragraph's Range for each metric substring in\n# left-to-right order, and set Font.Bold / Font.Color on the found hit.\n# Re-fetching $p.Range before every Find call keeps the search scoped to\n# the (growing) paragraph while still finding the next literal match.\n\n$d = $word.ActiveDocument\n\n# wdColor value for RRGGBB \"2C3E50\" is stored BGR-order: 0x00503E2C\n$highlightColor = 5258796\n\nfunction Set-MetricBold {\n    param([string]$LeadingText, [string[]]$Metrics)\n\n    $target = $null\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Text.TrimStart(\"`r\", \"`n\", \" \", [char]0x2022) -like \"$LeadingText*\") {\n            $target = $p\n            break\n        }\n    }\n    if ($null -eq $target) {\n        return\n    }\n\n    foreach ($metric in $Metrics) {\n        $rng = $target.Range\n        $find = $rng.Find\n        $find.Text = $metric\n        $find.MatchCase = $true\n        $find.MatchWildcards = $false\n        $find.Execute()\n        if ($find.Found) {\n            $rng.Font.Bold = $true\n            $rng.Font.Color = $highlightColor\n        }\n    }\n}\n\n$plusMinus = [char]0x00B1\n\nSet-MetricBold \"Discovered systematic race coding errors\" @(\"23%\", \"64%\")\nSet-MetricBold \"Utilized advanced sampling methods\" @($plusMinus + \"4.2%\", $plusMinus + \"2.1%\", \"71%\", \"87%\")\nSet-MetricBold \"Trigonometric algorithm for boundary estimation\" @(\"73.5%\", \"`$4.7M\")\nSet-MetricBold \"Built real-time FEC analysis systems\" @(\"`$2\")\nSet-MetricBold \"Algorithmic innovation: Pioneered trigonometric\" @(\"73.5%\")\nSet-MetricBold \"`$4.7M savings enabled\" @(\"`$4.7M\")\nSet-MetricBold \"178% accuracy improvement\" @(\"178%\")\n"}
